# --------------------------------------------------------------------------
# Refresh the cryptocurrency Price (col D) and Volume/1h (col E) figures on
# the active worksheet, as produced by the scheduled GitHub Actions scraper
# run referenced in the commit message. Every row's data is textual in the
# workbook (Price contains locale-formatted numbers such as "26.096.78" or
# leading-zero decimals such as "0.5450"), so each Price cell is explicitly
# formatted as Text before its value is written, preventing Excel from
# reinterpreting it as a number/date and silently dropping trailing zeros.
# --------------------------------------------------------------------------

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @(
    @{ Row = 2; D = '26.096.78'; E = '  -0.28%  ' },
    @{ Row = 3; D = '1.651.11'; E = '  -0.47%  ' },
    @{ Row = 4; E = '  -0.38%  ' },
    @{ Row = 5; D = '218.22'; E = '  +0.26%  ' },
    @{ Row = 6; D = '0.5294'; E = '  +1.48%  ' },
    @{ Row = 8; E = '  -2.06%  ' },
    @{ Row = 9; D = '0.06307'; E = '  -0.11%  ' },
    @{ Row = 10; E = '  -3.15%  ' },
    @{ Row = 11; D = '0.07749'; E = '  +0.34%  ' },
    @{ Row = 12; E = '  +1.00%  ' },
    @{ Row = 13; D = '1.652.09'; E = '  -0.50%  ' },
    @{ Row = 14; D = '0.5450'; E = '  +0.07%  ' },
    @{ Row = 15; D = '0.0₅8112'; E = '  -1.28%  ' },
    @{ Row = 16; D = '65.16'; E = '  +0.53%  ' },
    @{ Row = 17; D = '26.105.58'; E = '  -0.42%  ' },
    @{ Row = 18; D = '1.002'; E = '  -0.37%  ' },
    @{ Row = 19; D = '4.548'; E = '  -2.35%  ' },
    @{ Row = 20; D = '193.31'; E = '  +0.23%  ' },
    @{ Row = 21; E = '  -0.93%  ' },
    @{ Row = 22; E = '  -1.57%  ' },
    @{ Row = 23; E = '  -0.48%  ' },
    @{ Row = 24; D = '140.03'; E = '  +0.92%  ' },
    @{ Row = 25; D = '0.1240'; E = '  +0.42%  ' },
    @{ Row = 26; D = '7.240'; E = '  +0.58%  ' },
    @{ Row = 27; D = '16.13'; E = '  -0.04%  ' },
    @{ Row = 28; D = '1.436'; E = '  +1.45%  ' },
    @{ Row = 29; D = '0.05898'; E = '  -1.59%  ' },
    @{ Row = 30; D = '1.278'; E = '  -0.34%  ' },
    @{ Row = 31; E = '  -2.59%  ' },
    @{ Row = 32; D = '3.232'; E = '  -2.81%  ' },
    @{ Row = 33; D = '1.550'; E = '  -5.64%  ' },
    @{ Row = 34; D = '2.412'; E = '  -0.08%  ' },
    @{ Row = 35; D = '0.9434'; E = '  -3.63%  ' },
    @{ Row = 36; D = '2.757'; E = '  -1.02%  ' },
    @{ Row = 37; D = '0.5646'; E = '  -4.46%  ' },
    @{ Row = 38; D = '0.01602'; E = '  +0.96%  ' },
    @{ Row = 39; D = '5.836'; E = '  -1.80%  ' },
    @{ Row = 40; D = '0.8427'; E = '  -2.57%  ' },
    @{ Row = 41; D = '1.002'; E = '  -0.20%  ' },
    @{ Row = 42; D = '100.73'; E = '  +1.14%  ' },
    @{ Row = 43; D = '1.007.12'; E = '  -2.84%  ' },
    @{ Row = 44; D = '1.798.17'; E = '  -0.17%  ' },
    @{ Row = 45; D = '56.80'; E = '  -0.37%  ' },
    @{ Row = 46; E = '  -0.25%  ' },
    @{ Row = 47; E = '  +0.23%  ' },
    @{ Row = 48; D = '0.4286'; E = '  +1.26%  ' },
    @{ Row = 49; D = '1.472'; E = '  +1.14%  ' },
    @{ Row = 50; E = '  -0.61%  ' },
    @{ Row = 51; D = '7.776'; E = '  -3.69%  ' }
)

foreach ($u in $updates) {
    $rowNum = $u.Row
    if ($u.ContainsKey('D')) {
        $dCell = $ws.Range("D" + $rowNum)
        $dCell.NumberFormat = "@"
        $dCell.Value = $u.D
    }
    if ($u.ContainsKey('E')) {
        $ws.Range("E" + $rowNum).Value = $u.E
    }
}
